# Updated cryptos list on Wed Dec 20 04:59:02 UTC 2023 with GitHub Actions
#
# Refreshes the price/volume snapshot on the active sheet. Numeric-looking
# price strings are entered with a leading apostrophe so Excel keeps them
# as literal text (matching the sheet's existing "123.45" style strings)
# instead of silently re-parsing/rounding them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "42.745.93";  $ws.Range("E2").Value  = "  -0.55%  "
$ws.Range("D3").Value  = "2.215.37";   $ws.Range("E3").Value  = "  -1.06%  "
                                        $ws.Range("E4").Value  = "  -0.26%  "
$ws.Range("D5").Value  = "'253.54";    $ws.Range("E5").Value  = "  +2.76%  "
$ws.Range("D6").Value  = "'0.611";     $ws.Range("E6").Value  = "  -1.30%  "
$ws.Range("D7").Value  = "'75.55";     $ws.Range("E7").Value  = "  -0.67%  "
                                        $ws.Range("E8").Value  = "  -0.03%  "
$ws.Range("D9").Value  = "'0.589";     $ws.Range("E9").Value  = "  -4.41%  "
$ws.Range("D10").Value = "'40.98";     $ws.Range("E10").Value = "  +0.12%  "
                                        $ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "'6.89";      $ws.Range("E12").Value = "  -1.25%  "
                                        $ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "2.546.02";   $ws.Range("E14").Value = "  -0.31%  "
                                        $ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "2.219.83";   $ws.Range("E16").Value = "  -0.79%  "
                                        $ws.Range("E17").Value = "  -3.88%  "
$ws.Range("D18").Value = "42.655.30"
                                        $ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").Value = "'71.14";     $ws.Range("E20").Value = "  +0.01%  "
                                        $ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").Value = "'2.19";      $ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "'228.99";    $ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "'9.50";      $ws.Range("E24").Value = "  -8.82%  "
                                        $ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'10.60";     $ws.Range("E26").Value = "  -2.96%  "
                                        $ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'39.26";     $ws.Range("E28").Value = "  +3.76%  "

# Rows 29/30 swapped rank: Toncoin now outranks PancakeSwap.
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.17"
$ws.Range("E30").Value = "  -3.61%  "

$ws.Range("D31").Value = "'173.40";    $ws.Range("E31").Value = "  +0.03%  "
                                        $ws.Range("E32").Value = "  -0.34%  "
                                        $ws.Range("E33").Value = "  +4.43%  "
                                        $ws.Range("E34").Value = "  -3.37%  "
                                        $ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").Value = "'0.0344";    $ws.Range("E37").Value = "  +4.81%  "
$ws.Range("D38").Value = "'4.29";      $ws.Range("E38").Value = "  -1.17%  "
                                        $ws.Range("E39").Value = "  -5.71%  "
                                        $ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").Value = "'2.72";      $ws.Range("E41").Value = "  +17.73%  "
$ws.Range("D42").Value = "'5.27";      $ws.Range("E42").Value = "  -5.44%  "
$ws.Range("D43").Value = "'59.91";     $ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").Value = "'0.197";     $ws.Range("E44").Value = "  -1.86%  "
$ws.Range("D45").Value = "'103.08";    $ws.Range("E45").Value = "  -2.24%  "
$ws.Range("D46").Value = "'8.36";      $ws.Range("E46").Value = "  -2.85%  "
                                        $ws.Range("E47").Value = "  -0.75%  "
                                        $ws.Range("E48").Value = "  +3.69%  "
                                        $ws.Range("E49").Value = "  -0.62%  "

# Row 51: HuobiToken dropped out of the top 50, replaced by RocketPoolETH.
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.440.74"
$ws.Range("E51").Value = "  +0.02%  "
